$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: correct the "raw password" column (C) so rows that should use the
# capitalized password reflect the real value used by the automation run.
$ws.Cells.Item(3, 3).Value = "Password1"
$ws.Cells.Item(5, 3).Value = "Password1"
$ws.Cells.Item(6, 3).Value = "Password1"

# Step 2: add the new "expected" outcome column (D) for each existing row.
$ws.Cells.Item(1, 4).Value = "expected"
$ws.Cells.Item(2, 4).Value = "fail"
$ws.Cells.Item(3, 4).Value = "fail"
$ws.Cells.Item(4, 4).Value = "fail"
$ws.Cells.Item(5, 4).Value = "success"
$ws.Cells.Item(6, 4).Value = "success"

# Step 3: row 5 now reuses the automation_user / automation password combo
# (previously a numeric placeholder row) and row 6 becomes the administrator
# success case, replacing the old petrolink123 / "a" rows.
$ws.Cells.Item(5, 1).Value = "automation_user"
$ws.Cells.Item(5, 2).Value = "p4y+y39Ir5Oy1MY8jPt0uQ=="
$ws.Cells.Item(6, 1).Value = "administrator"
$ws.Cells.Item(6, 2).Value = "p4y+y39Ir5Oy1MY8jPt0uQ=="

# The table shrank from 7 to 6 rows, so clear out the old trailing row.
$ws.Range("A7:C7").ClearContents()

# Update the active selection to match the target state
$ws.Range("G3").Select()
